# "add the NA's under duplicate_image_filename"
# Column E (header "duplicate_image_filename") gets the literal string
# "NA" added for every data row that currently has content (rows 2-21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2:E21").Value = "NA"
